$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aggiornamenti_29_08_2017")

# --- Step 1: copy existing rows' formatting down into the new rows so the
#     new cells inherit the same cellXfs (border/wrap/date-format) as their
#     closest analogue already in the table; row height is fixed up after. ---
$ws.Range("A84:D84").Copy()
$ws.Range("A85:D85").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A86:D86").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A87:D87").PasteSpecial(-4122) | Out-Null
$ws.Range("A72:D72").Copy()
$ws.Range("A88:D88").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A89:D89").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A90:D90").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A91:D91").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A92:D92").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A93:D93").PasteSpecial(-4122) | Out-Null
$ws.Range("A47:D47").Copy()
$ws.Range("A94:D94").PasteSpecial(-4122) | Out-Null
$ws.Range("A47:D47").Copy()
$ws.Range("A95:D95").PasteSpecial(-4122) | Out-Null
$ws.Range("A84:D84").Copy()
$ws.Range("A96:D96").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 1b: a few individual cells need a different donor style than the
#     rest of their row (kept as single-cell copies to match exactly). ---
$ws.Range("C57").Copy()
$ws.Range("C90").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 2: set shared-string cell values in the exact order that makes
#     new entries land at the same sharedStrings.xml indices as the target. ---
$ws.Range("C85").Value = "tabella 3 archivio comuni 20171005.xlsx"
$ws.Range("B85").Value = "Tabella 3 comuni"
$ws.Range("D85").Value = "La denominazione per MONTEBELLO IONICO diventa MONTEBELLO JONICO `nIl codice catastale del comune CASALI DEL MANCO è impostato a M385 (prima era N.D, i.e. Non Disponibile)`nIl codice catastale per OLGIATE CALCO (due record) diventa G027 (prima era G026) in accordo a quanto presente nella banca dati di AE`n"
$ws.Range("C86").Value = "errori_anpr_05102017.xlsx"
$ws.Range("D86").Value = "Inserimento codice di errore EN407, EN427`n"
$ws.Range("C87").Value = "errori_anpr_11102017.xlsx"
$ws.Range("D87").Value = "Eliminato codice di errore EN375 dal servizio 2003`n"
$ws.Range("B88").Value = "tipoDatiControllo"
$ws.Range("D88").Value = "Aggiornato il file vocabolario5000mutazione.xsd per documentare l'utilizzo dei dati di controllo "
$ws.Range("C89").Value = "errori_anpr_16102017.xlsx"
$ws.Range("D89").Value = "Il codice di errore EC039 è attivo anche per il subentro`n"
$ws.Range("D90").Value = "In Allegato XML S001- Subentro.xls specificate regole/condizioni per convivenza e responsabile convivenza"
$ws.Range("C91").Value = "errori_anpr_18102017.xlsx"
$ws.Range("D91").Value = "Inserimento codice di errore EN436"
$ws.Range("C92").Value = "errori_anpr_19102017.xlsx"
$ws.Range("D92").Value = "Inserimento codice di errore EN411"
$ws.Range("B93").Value = "Invio file di Subentro"
$ws.Range("D93").Value = "Aggiornate istruzioni per la predisposizione del file AIRE con AnagAire 6.0.3"
$ws.Range("C94").Value = "Sito WEB di ANPR e specifiche di integrazione.pdf"
$ws.Range("C95").Value = "Allegato 2 - Elenco funzioni WEB19102017.xlsx"
$ws.Range("C96").Value = "Allegato 7 - Utilizzo WS ANPR totale 19102017.xlsx"
$ws.Range("D94").Value = "Inserita la descrizione delle seguenti funzioni:`nRegistrazione/Eliminazione dati`nRegistrazione/Rettifiche`nConsultazione/Consultazione AE"
$ws.Range("D96").Value = "Inserita operazione anagrafica 4002"
$ws.Range("B94").Value = "Documentazione sito WEB"

# --- Step 3: fill in cells that reuse pre-existing shared strings. ---
$ws.Range("B86").Value = "lista errori"
$ws.Range("B87").Value = "lista errori"
$ws.Range("C88").Value = "Tracciati XSD e WSDL - rar"
$ws.Range("B89").Value = "lista errori"
$ws.Range("B90").Value = "specifiche xlsx"
$ws.Range("C90").Value = "Specifiche di interfaccia dei servizi di ANPR per i comuni – documentazione tecnica - rar"
$ws.Range("B91").Value = "lista errori"
$ws.Range("B92").Value = "lista errori"
$ws.Range("C93").Value = "Invio file di Subentro.pdf"
$ws.Range("B95").Value = "Documentazione sito WEB"
$ws.Range("D95").Value = "Inserita la descrizione delle seguenti funzioni:`nRegistrazione/Eliminazione dati`nRegistrazione/Rettifiche`nConsultazione/Consultazione AE"
$ws.Range("B96").Value = "Documentazione sito WEB"

# --- Step 4: date column A + row heights ---
$ws.Range("A85").Value = 42865
$ws.Range("A86").Value = 42865
$ws.Range("A87").Value = 43019
$ws.Range("A88").Value = 43019
$ws.Range("A89").Value = 43024
$ws.Range("A90").Value = 43025
$ws.Range("A91").Value = 43026
$ws.Range("A92").Value = 43027
$ws.Range("A93").Value = 43027
$ws.Range("A94").Value = 43027
$ws.Range("A95").Value = 43027
$ws.Range("A96").Value = 43027
$ws.Rows.Item(85).RowHeight = 93
$ws.Rows.Item(86).RowHeight = 30
$ws.Rows.Item(87).RowHeight = 30
$ws.Rows.Item(88).RowHeight = 30
$ws.Rows.Item(89).RowHeight = 30
$ws.Rows.Item(90).RowHeight = 30
$ws.Rows.Item(94).RowHeight = 60
$ws.Rows.Item(95).RowHeight = 60

# --- Step 5: sheet view / selection state, matching the post-edit XML ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 86
$ws.Range("B95:B96").Select() | Out-Null
